$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New taxonomy-style rows to append below the existing label/content pairs.
# Copy the header row's formatting (style index used by row 1) onto the
# new rows before filling in their values, so the new cells are styled
# the same way as the "label"/"content" header instead of the plain
# body style used by rows 2-20.
$ws.Range("A1:B1").Copy()
$ws.Range("A21:B27").PasteSpecial(-4122) # xlPasteFormats

$rows = @(
    @("Class", "Submissions"),
    @("Group", "Joint Locks"),
    @("Family", "Wrist Locks"),
    @("SubFamily", "Flexion"),
    @("Genus", "Gooseneck"),
    @("Species", "Gooseneck from Guard"),
    @("Variety", "Gooseneck form Half Guard")
)

$startRow = 21
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
